$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test case row 18 (mirrors existing BWP bootstrap rows)
$ws.Range("A18").Value = "MissingReqFields"
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 950
$ws.Range("D18").Value = 1.5
$ws.Range("O18").Value = "udf data 9"

# Move selection to reflect the new last-used cell
$ws.Range("E18").Select()
